# Generate Report for Handback
#
# The CI run that produced this handback-status.xlsx was re-run: the two
# e2e fixture files were renamed (new GUIDs) and the handback/handoff
# timestamps + xlf hash advanced. This script mirrors that regeneration
# across all three sheets (Overview, zh-cn, de-de) plus the hyperlink
# display text that mirrors the file-name cells.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---- old -> new identifiers -------------------------------------------------
$oldGuid1 = "07eaff7d-56da-4823-9679-719e2d7cef8a"
$newGuid1 = "fdf571cb-597c-45aa-a598-20719cce4dc6"

$oldGuid2 = "be1543a8-66a4-4809-9562-42f4569a305a"
$newGuid2 = "ffffe5b7fcbf-c159-4dfa-97f2-708c9665f52b"

$newXlfHash = "62841f88e27d8f567c8f53ea8e14d619c0249069"

$newFile1 = $newGuid1 + ".md"
$newFile2 = $newGuid2 + ".md"

$newXlfZhCn = $newGuid1 + "." + $newXlfHash + ".zh-cn.xlf"
$newXlfDeDe = $newGuid1 + "." + $newXlfHash + ".de-de.xlf"

# ---- Overview sheet -----------------------------------------------------
$ws1.Range("A2").Value = $newFile1
$ws1.Range("B2").Value = "e2e\" + $newFile1
$ws1.Range("G2").Value = "2016-09-06 09:26:58"

$ws1.Range("A3").Value = $newFile2
$ws1.Range("B3").Value = "e2e\" + $newFile2
$ws1.Range("G3").Value = "2016-09-06 09:26:58"

# ---- zh-cn sheet ----------------------------------------------------------
$ws2.Range("A2").Value = $newFile1
$ws2.Range("G2").Value = $newXlfZhCn
$ws2.Range("H2").Value = "2016-09-06 09:26:53"
$ws2.Range("I2").Value = $newFile1
$ws2.Range("J2").Value = $newXlfZhCn
$ws2.Range("K2").Value = "2016-09-06 09:27:25"

$ws2.Range("A3").Value = $newFile2
$ws2.Range("G3").Value = $newXlfZhCn
$ws2.Range("H3").Value = "2016-09-06 09:26:53"
$ws2.Range("I3").Value = $newFile2
$ws2.Range("J3").Value = $newXlfZhCn
$ws2.Range("K3").Value = "2016-09-06 09:27:25"

# ---- de-de sheet ------------------------------------------------------------
$ws3.Range("A2").Value = $newFile1
$ws3.Range("G2").Value = $newXlfDeDe
$ws3.Range("H2").Value = "2016-09-06 09:26:58"
$ws3.Range("I2").Value = $newFile1
$ws3.Range("J2").Value = $newXlfDeDe
$ws3.Range("K2").Value = "2016-09-06 09:27:32"

$ws3.Range("A3").Value = $newFile2
$ws3.Range("G3").Value = $newXlfDeDe
$ws3.Range("H3").Value = "2016-09-06 09:26:58"
$ws3.Range("I3").Value = $newFile2
$ws3.Range("J3").Value = $newXlfDeDe
$ws3.Range("K3").Value = "2016-09-06 09:27:32"

# ---- hyperlink display text (independent of the cell text) ------------------
foreach ($ws in @($ws1, $ws2, $ws3)) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.TextToDisplay -eq ("e2e\" + $oldGuid1 + ".md")) {
            $hl.TextToDisplay = "e2e\" + $newFile1
        } elseif ($hl.TextToDisplay -eq ("e2e\" + $oldGuid2 + ".md")) {
            $hl.TextToDisplay = "e2e\" + $newFile2
        } elseif ($hl.TextToDisplay -eq ($oldGuid1 + ".md")) {
            $hl.TextToDisplay = $newFile1
        } elseif ($hl.TextToDisplay -eq ($oldGuid2 + ".md")) {
            $hl.TextToDisplay = $newFile2
        }
    }
}
